$d = $word.ActiveDocument

# The last paragraph in the document is the final bullet item in the
# "extra credit" list ("... print function for visualization").
# We append a new bullet item after it, inheriting the same
# ListParagraph style / numbering (numId 1, ilvl 0) that the
# surrounding list items already use.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$r = $newPara.Range
$r.Collapse(0)

$dash = [char]0x2013

$r.InsertAfter("Maven " + $dash + " u")
$r.Collapse(0)
$r.InsertAfter("s")
$r.Collapse(0)
$r.InsertAfter("ed for JFreeApache + mini report (learning Maven)")
$r.Collapse(0)
$r.InsertAfter(" + connected Maven with GitHub action")
$r.Collapse(0)
$r.InsertAfter("s")
$r.Collapse(0)
$r.InsertAfter(" for automatic building")
